# Updated cryptos list on Mon Feb  5 04:44:12 UTC 2024 with GitHub Actions
#
# Refreshes the crypto market snapshot on the active worksheet: each row's
# Price (column D) and Volume(1h) (column E) are updated to the latest
# scraped figures, and three rows (44-46) shift rank -- EnergySwap drops out
# of its slot and FraxShare/ApeXProtocol/BitcoinSV move in -- mirroring the
# upstream coinranking.com re-ranking for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'42.734.39"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "'2.294.48"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'303.58"
$ws.Range("E5").Value = "  +1.30%  "

# Row 6
$ws.Range("D6").Value = "'96.58"
$ws.Range("E6").Value = "  -1.01%  "

# Row 7
$ws.Range("E7").Value = "  -2.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -2.01%  "

# Row 10
$ws.Range("D10").Value = "'35.00"
$ws.Range("E10").Value = "  -2.20%  "

# Row 11
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -0.86%  "

# Row 12
$ws.Range("D12").Value = "'18.71"
$ws.Range("E12").Value = "  +5.70%  "

# Row 13
$ws.Range("E13").Value = "  +2.19%  "

# Row 14
$ws.Range("E14").Value = "  +0.87%  "

# Row 15
$ws.Range("D15").Value = "'2.652.04"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16
$ws.Range("D16").Value = "'2.293.72"
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("E17").Value = "  -0.41%  "

# Row 18
$ws.Range("D18").Value = "'42.661.89"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19
$ws.Range("D19").Value = "'12.80"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0894"
$ws.Range("E20").Value = "  -1.51%  "

# Row 21
$ws.Range("D21").Value = "'5.99"
$ws.Range("E21").Value = "  -1.76%  "

# Row 22
$ws.Range("D22").Value = "'67.21"
$ws.Range("E22").Value = "  -1.29%  "

# Row 23
$ws.Range("D23").Value = "'235.90"
$ws.Range("E23").Value = "  -2.29%  "

# Row 24
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  -1.74%  "

# Row 27
$ws.Range("D27").Value = "'24.98"
$ws.Range("E27").Value = "  -0.57%  "

# Row 28
$ws.Range("D28").Value = "'167.49"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("E29").Value = "  +1.06%  "

# Row 30
$ws.Range("E30").Value = "  -0.54%  "

# Row 31
$ws.Range("D31").Value = "'32.90"
$ws.Range("E31").Value = "  -0.15%  "

# Row 32
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("D33").Value = "'17.90"
$ws.Range("E33").Value = "  +1.70%  "

# Row 34
$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = "  -0.94%  "

# Row 35
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "  -4.89%  "

# Row 36
$ws.Range("E36").Value = "  -1.87%  "

# Row 37
$ws.Range("E37").Value = "  -0.42%  "

# Row 38
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("E39").Value = "  -0.99%  "

# Row 40
$ws.Range("E40").Value = "  -1.17%  "

# Row 41
$ws.Range("E41").Value = "  -2.60%  "

# Row 42
$ws.Range("D42").Value = "'1.993.56"
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("D43").Value = "'0.0279"
$ws.Range("E43").Value = "  -2.44%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'10.11"
$ws.Range("E44").Value = "  -0.21%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'18.07"
$ws.Range("E45").Value = "  +4.45%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = "  -1.29%  "

# Row 47
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  -0.98%  "

# Row 49
$ws.Range("D49").Value = "'53.61"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "'2.518.81"
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'70.82"
$ws.Range("E51").Value = "  -2.08%  "
